$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F4").Value = -3
$ws.Range("F7").Value = 4
$ws.Range("F9").Value = -4
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = -4
$ws.Range("F16").Value = -7
$ws.Range("F23").Value = -6
$ws.Range("F31").Value = -1
$ws.Range("F32").Value = 3
